$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country label swaps (ranking shuffled by the refreshed case counts) ---
# Barein / Hungria swap around row 63-64
$ws.Range("A63").Value = "Hungria"
$ws.Range("A64").Value = "Barein"

# Azerbaiyan / Uzbekistan swap around row 71-72
$ws.Range("A71").Value = "Uzbekistan"
$ws.Range("A72").Value = "Azerbaiyan"

# Afganistan / Bulgaria swap around row 84-85
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("A85").Value = "Afganistan"

# --- Updated timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 07:22"

# --- Updated statistics (row 4: Estados Unidos) ---
$ws.Range("B4").Value = 614246
$ws.Range("C4").Value = 360
$ws.Range("E4").Value = 549362

# --- Updated statistics (row 16) ---
$ws.Range("D16").Value = 14700
$ws.Range("E16").Value = 10062

# --- Updated statistics (row 63, now Hungria) ---
$ws.Range("B63").Value = 1579
$ws.Range("C63").Value = 67
$ws.Range("D63").Value = 192
$ws.Range("E63").Value = 1253
$ws.Range("F63").Value = 58
$ws.Range("G63").Value = 12
$ws.Range("H63").Value = 134

# --- Updated statistics (row 64, now Barein) ---
$ws.Range("B64").Value = 1528
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 645
$ws.Range("E64").Value = 876
$ws.Range("F64").Value = 3
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 7

# --- Updated statistics (row 71, now Uzbekistan) ---
$ws.Range("B71").Value = 1214
$ws.Range("C71").Value = 49
$ws.Range("D71").Value = 99
$ws.Range("E71").Value = 1111
$ws.Range("F71").Value = 8
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 4

# --- Updated statistics (row 72, now Azerbaiyan) ---
$ws.Range("B72").Value = 1197
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 351
$ws.Range("E72").Value = 833
$ws.Range("F72").Value = 25
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 13

# --- Updated statistics (row 84, now Bulgaria) ---
$ws.Range("B84").Value = 735
$ws.Range("C84").Value = 22
$ws.Range("D84").Value = 105
$ws.Range("E84").Value = 594
$ws.Range("F84").Value = 29
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 36

# --- Updated statistics (row 85, now Afganistan) ---
$ws.Range("B85").Value = 714
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 40
$ws.Range("E85").Value = 651
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 23

# --- Updated statistics (row 99: Honduras) ---
$ws.Range("D99").Value = 9
$ws.Range("E99").Value = 379
